{"js": "// Replace word-problem text, number-sentence blanks, and answer blanks\n// for each of the 7 practice problems in the document.\n\nconst replacements = [\n  {\n    oldQ: \"Lily has eight crayons. Sam gives her five more. How many crayons does Lily have in total?\",\n    newQ: \"Barnaby has seven bouncy balls. He buys three more. How many bouncy balls does Barnaby have in total?\",\n    oldA: \"Answer: Lily has ____________ crayons in total.\",\n    newA: \"Answer: Barnaby has ________________________________________ bouncy balls in total.\"\n  },\n  {\n    oldQ: \"David had twelve toy cars. He gave three of them to his friend. How many toy cars does David have left?\",\n    newQ: \"Penelope has twelve rainbow erasers. She gives four to her friend Herbert. How many rainbow erasers does Penelope have left?\",\n    oldA: \"Answer: David has ____________ toy cars left.\",\n    newA: \"Answer: Penelope has ________________________________________ rainbow erasers left.\"\n  },\n  {\n    oldQ: \"Sarah baked seven cookies. Her brother ate one cookie. How many cookies are left?\",\n    newQ: \"Freddy found eight shiny rocks. He found another five shiny rocks. How many shiny rocks does Freddy have now?\",\n    oldA: \"Answer: There are ____________ cookies left.\",\n    newA: \"Answer: Freddy now has ________________________________________ shiny rocks.\"\n  },\n  {\n    oldQ: \"Emily found nine shells at the beach. Her dad found six shells. How many shells did they find altogether?\",\n    newQ: \"Gertrude baked fifteen cookies. She ate three of them. How many cookies are left?\",\n    oldA: \"Answer: They found ____________ shells altogether.\",\n    newA: \"Answer: There are ________________________________________ cookies left.\"\n  },\n  {\n    oldQ: \"Jack has four apples and buys eleven more. How many apples does Jack have now?\",\n    newQ: \"Kevin saw six purple frogs. Then two more purple frogs hopped along. How many purple frogs did Kevin see in all?\",\n    oldA: \"Answer: Jack now has ____________ apples.\",\n    newA: \"Answer: Kevin saw ________________________________________ purple frogs in all.\"\n  },\n  {\n    oldQ: \"Chloe had thirteen balloons, and four of them popped. How many balloons does Chloe have left?\",\n    newQ: \"Esmeralda has nine sparkly stickers. She uses two of them on her notebook. How many stickers does Esmeralda have remaining?\",\n    oldA: \"Answer: Chloe has ____________ balloons left.\",\n    newA: \"Answer: Esmeralda has ________________________________________ stickers remaining.\"\n  },\n  {\n    oldQ: \"Ben has ten stickers, then he gets five more stickers for his birthday. How many stickers does Ben have in all?\",\n    newQ: \"Humphrey collected five acorns. Then, a squirrel gave him one more acorn. How many acorns does Humphrey have?\",\n    oldA: \"Answer: Ben has ____________ stickers in all.\",\n    newA: \"Answer: Humphrey has ________________________________________ acorns.\"\n  }\n];\n\nconst oldNumberSentence = \"Number sentence: __________________________________________________\";\nconst newNumberSentence = \"Number sentence: ______________________________________________________________________\";\n\nconst body = context.document.body;\n\n// Replace each unique question/answer string in document order.\nfor (const r of replacements) {\n  const qRanges = body.search(r.oldQ, { matchCase: true, matchWholeWord: false });\n  qRanges.load(\"items\");\n  await context.sync();\n  for (const rng of qRanges.items) {\n    rng.insertText(r.newQ, \"Replace\");\n  }\n  await context.sync();\n\n  const aRanges = body.search(r.oldA, { matchCase: true, matchWholeWord: false });\n  aRanges.load(\"items\");\n  await context.sync();\n  for (const rng of aRanges.items) {\n    rng.insertText(r.newA, \"Replace\");\n  }\n  await context.sync();\n}\n\n// Replace all 7 \"Number sentence:\" blank lines uniformly.\nconst nsRanges = body.search(oldNumberSentence, { matchCase: true, matchWholeWord: false });\nnsRanges.load(\"items\");\nawait context.sync();\nfor (const rng of nsRanges.items) {\n  rng.insertText(newNumberSentence, \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Replace word-problem text, number-sentence blanks, and answer blanks\n# for each of the 7 practice problems in the document.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{\n        OldQ = \"Lily has eight crayons. Sam gives her five more. How many crayons does Lily have in total?\"\n        NewQ = \"Barnaby has seven bouncy balls. He buys three more. How many bouncy balls does Barnaby have in total?\"\n        OldA = \"Answer: Lily has ____________ crayons in total.\"\n        NewA = \"Answer: Barnaby has ________________________________________ bouncy balls in total.\"\n    },\n    @{\n        OldQ = \"David had twelve toy cars. He gave three of them to his friend. How many toy cars does David have left?\"\n        NewQ = \"Penelope has twelve rainbow erasers. She gives four to her friend Herbert. How many rainbow erasers does Penelope have left?\"\n        OldA = \"Answer: David has ____________ toy cars left.\"\n        NewA = \"Answer: Penelope has ________________________________________ rainbow erasers left.\"\n    },\n    @{\n        OldQ = \"Sarah baked seven cookies. Her brother ate one cookie. How many cookies are left?\"\n        NewQ = \"Freddy found eight shiny rocks. He found another five shiny rocks. How many shiny rocks does Freddy have now?\"\n        OldA = \"Answer: There are ____________ cookies left.\"\n        NewA = \"Answer: Freddy now has ________________________________________ shiny rocks.\"\n    },\n    @{\n        OldQ = \"Emily found nine shells at the beach. Her dad found six shells. How many shells did they find altogether?\"\n        NewQ = \"Gertrude baked fifteen cookies. She ate three of them. How many cookies are left?\"\n        OldA = \"Answer: They found ____________ shells altogether.\"\n        NewA = \"Answer: There are ________________________________________ cookies left.\"\n    },\n    @{\n        OldQ = \"Jack has four apples and buys eleven more. How many apples does Jack have now?\"\n        NewQ = \"Kevin saw six purple frogs. Then two more purple frogs hopped along. How many purple frogs did Kevin see in all?\"\n        OldA = \"Answer: Jack now has ____________ apples.\"\n        NewA = \"Answer: Kevin saw ________________________________________ purple frogs in all.\"\n    },\n    @{\n        OldQ = \"Chloe had thirteen balloons, and four of them popped. How many balloons does Chloe have left?\"\n        NewQ = \"Esmeralda has nine sparkly stickers. She uses two of them on her notebook. How many stickers does Esmeralda have remaining?\"\n        OldA = \"Answer: Chloe has ____________ balloons left.\"\n        NewA = \"Answer: Esmeralda has ________________________________________ stickers remaining.\"\n    },\n    @{\n        OldQ = \"Ben has ten stickers, then he gets five more stickers for his birthday. How many stickers does Ben have in all?\"\n        NewQ = \"Humphrey collected five acorns. Then, a squirrel gave him one more acorn. How many acorns does Humphrey have?\"\n        OldA = \"Answer: Ben has ____________ stickers in all.\"\n        NewA = \"Answer: Humphrey has ________________________________________ acorns.\"\n    }\n)\n\n$oldNumberSentence = \"Number sentence: __________________________________________________\"\n$newNumberSentence = \"Number sentence: ______________________________________________________________________\"\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.Execute($r.OldQ, $false, $false, $false, $false, $false, $true, 1, $false, $r.NewQ, 2) | Out-Null\n\n    $find2 = $d.Content.Find\n    $find2.Execute($r.OldA, $false, $false, $false, $false, $false, $true, 1, $false, $r.NewA, 2) | Out-Null\n}\n\n# Replace all 7 \"Number sentence:\" blank lines uniformly.\n$findNs = $d.Content.Find\n$findNs.Execute($oldNumberSentence, $false, $false, $false, $false, $false, $true, 1, $false, $newNumberSentence, 2) | Out-Null\n"}
